$d = $word.ActiveDocument

# The "References" section lists each reference's name in bold (e.g. the
# existing "Bishal Bhandari" paragraph), followed by their title/employer in
# regular weight. The "Prem Acharya" reference paragraph was missing that
# bold styling on both its run and its paragraph mark - add it, matching the
# pattern used elsewhere in the document (w:b + w:bCs on both the paragraph
# mark's rPr and the run's rPr).

$range = $d.Content
$found = $range.Find.Execute("Prem Acharya", $false, $false, $false, $false,
                              $false, $true, 1, $false, "", 0)

if ($found -and $range.Find.Found) {
    $para = $range.Paragraphs(1)
    $paraRange = $para.Range

    $xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
        '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
        '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
        '<pkg:xmlData>' +
        '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml">' +
        '<w:body>' +
        '<w:p w14:paraId="7525CDE4" w14:textId="38645D6D" w:rsidR="007130F5" w:rsidRDefault="007130F5" w:rsidP="007130F5">' +
          '<w:pPr>' +
            '<w:widowControl w:val="0"/>' +
            '<w:pBdr>' +
              '<w:top w:val="nil"/>' +
              '<w:left w:val="nil"/>' +
              '<w:bottom w:val="nil"/>' +
              '<w:right w:val="nil"/>' +
              '<w:between w:val="nil"/>' +
            '</w:pBdr>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Times" w:eastAsia="Times" w:hAnsi="Times" w:cs="Times"/>' +
              '<w:b/>' +
              '<w:bCs/>' +
              '<w:color w:val="000000"/>' +
              '<w:sz w:val="20"/>' +
              '<w:szCs w:val="20"/>' +
            '</w:rPr>' +
          '</w:pPr>' +
          '<w:r>' +
            '<w:rPr>' +
              '<w:rFonts w:ascii="Times" w:eastAsia="Times" w:hAnsi="Times" w:cs="Times"/>' +
              '<w:b/>' +
              '<w:bCs/>' +
              '<w:color w:val="000000"/>' +
              '<w:sz w:val="20"/>' +
              '<w:szCs w:val="20"/>' +
            '</w:rPr>' +
            '<w:t>Prem Acharya</w:t>' +
          '</w:r>' +
        '</w:p>' +
        '</w:body>' +
        '</w:document>' +
        '</pkg:xmlData></pkg:part></pkg:package>'

    $paraRange.InsertXML($xml)
}
